$d = $word.ActiveDocument

# 1) Update object hash identifiers in the stack trace text
$null = $d.Content.Find.Execute("EClassImpl@6e818345", $true, $false, $false, $false, $false, $true, 1, $false, "EClassImpl@582c1f8d", 2)
$null = $d.Content.Find.Execute("EClassImpl@25a5c7db", $true, $false, $false, $false, $false, $true, 1, $false, "EClassImpl@6f27a732", 2)

# 2) Update line numbers referenced in the test-suite frames
$null = $d.Content.Find.Execute("AbstractTemplatesTestSuite.java:518", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:536", 2)
$null = $d.Content.Find.Execute("AbstractTemplatesTestSuite.java:414", $true, $false, $false, $false, $false, $true, 1, $false, "AbstractTemplatesTestSuite.java:422", 2)

# 3) Update the generated reflective accessor number
$null = $d.Content.Find.Execute("GeneratedMethodAccessor5", $true, $false, $false, $false, $false, $true, 1, $false, "GeneratedMethodAccessor4", 2)

# 4) Replace the Maven/Tycho/Equinox launcher tail of the stack trace with the
#    Eclipse JDT JUnit runner tail (reflects running the test from within the IDE)
$oldBlock = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n" +
"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n" +
"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n" +
"`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
"`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
"`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n" +
"`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n" +
"`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n" +
"`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n" +
"`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
"`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
"`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n" +
"`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n" +
"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n" +
"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n" +
"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n" +
"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
"`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
"`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n" +
"`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n" +
"`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n" +
"`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)`n"

$newBlock = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" +
"`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" +
"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" +
"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" +
"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" +
"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"

$null = $d.Content.Find.Execute($oldBlock, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock, 2)

Write-Host "done"
